$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna5", "Epha4", "ECs", 1, 0.3333333333333333, 0.2457683333333333, 0.737305, 0.1447271191911903, 0.1575855905380038, 3, 1, 4.886733666666667, 14.660201, 0.4081653954827624, 0.4171266852711343, 1.201004388700556, 10.809039498305, 0.05907260184175309, 0.06573315502761172),
    @("ECs", "Efna5", "Epha4", "FAPs", 1, 0.3333333333333333, 0.2457683333333333, 0.737305, 0.1447271191911903, 0.1575855905380038, 3, 1, 6.292848333333333, 18.878545, 0.5256114009667484, 0.5371512231375235, 1.546582846802778, 13.919245621225, 0.07607022387596313, 0.08464729270633767),
    @("ECs", "Efna5", "Epha4", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2457683333333333, 0.737305, 0.1447271191911903, 0.1575855905380038, 1, 0.3333333333333333, 0.021228, 0.063684, 0.001773072896198643, 0.001812000792131494, 0.00521717018, 0.04695453162, 0.00025661173238281, 0.000285545214883372),
    @("ECs", "Efna5", "Epha4", "MuSCs", 1, 0.3333333333333333, 0.2457683333333333, 0.737305, 0.1447271191911903, 0.1575855905380038, 2, 1, 0.771625, 1.54325, 0.06445013065429045, 0.04391009079921059, 0.1896409902083333, 1.13784594125, 0.009327681741091283, 0.006919597589170966),
    @("FAPs", "Efna5", "Epha4", "ECs", 3, 1, 1.036691, 3.110073, 0.6104826439049008, 0.6647217777192627, 3, 1, 4.886733666666667, 14.660201, 0.4081653954827624, 0.4171266852711343, 5.066032811630334, 45.594295304673, 0.2491778897848063, 0.2772731917675718),
    @("FAPs", "Efna5", "Epha4", "FAPs", 3, 1, 1.036691, 3.110073, 0.6104826439049008, 0.6647217777192627, 3, 1, 6.292848333333333, 18.878545, 0.5256114009667484, 0.5371512231375235, 6.523739231531666, 58.713653083785, 0.3208766377287395, 0.357056115948051),
    @("FAPs", "Efna5", "Epha4", "Inflammatory-Mac", 3, 1, 1.036691, 3.110073, 0.6104826439049008, 0.6647217777192627, 1, 0.3333333333333333, 0.021228, 0.063684, 0.001773072896198643, 0.001812000792131494, 0.022006876548, 0.198061888932, 0.001082430229507467, 0.001204476387774359),
    @("FAPs", "Efna5", "Epha4", "MuSCs", 3, 1, 1.036691, 3.110073, 0.6104826439049008, 0.6647217777192627, 2, 1, 0.771625, 1.54325, 0.06445013065429045, 0.04391009079921059, 0.799936692875, 4.79962015725, 0.03934568616184753, 0.0291879936158655),
    @("MuSCs", "Efna5", "Epha4", "ECs", 2, 1, 0.4156905, 0.831381, 0.2447902369039089, 0.1776926317427335, 3, 1, 4.886733666666667, 14.660201, 0.4081653954827624, 0.4171266852711343, 2.0313687612635, 12.188212567581, 0.09991490385620307, 0.07412033847595074),
    @("MuSCs", "Efna5", "Epha4", "FAPs", 2, 1, 0.4156905, 0.831381, 0.2447902369039089, 0.1776926317427335, 3, 1, 6.292848333333333, 18.878545, 0.5256114009667484, 0.5371512231375235, 2.6158772701075, 15.695263620645, 0.1286645393620458, 0.09544781448313482),
    @("MuSCs", "Efna5", "Epha4", "Inflammatory-Mac", 2, 1, 0.4156905, 0.831381, 0.2447902369039089, 0.1776926317427335, 1, 0.3333333333333333, 0.021228, 0.063684, 0.001773072896198643, 0.001812000792131494, 0.008824277934, 0.052945667604, 0.0004340309343083656, 0.0003219791894737629),
    @("MuSCs", "Efna5", "Epha4", "MuSCs", 2, 1, 0.4156905, 0.831381, 0.2447902369039089, 0.1776926317427335, 2, 1, 0.771625, 1.54325, 0.06445013065429045, 0.04391009079921059, 0.3207571820625, 1.28302872825, 0.01577676275135164, 0.007802499594174117)
)

# Delete rows 14-16 (no longer needed since Resolving-Mac target cluster was removed)
$ws.Rows.Item(14).Resize(3).Delete() | Out-Null

# Write updated data (rows 2-13) reflecting the new TPM-based calculations
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
